# "Include Formula in output":
# Add a formula to C2 on Tabelle1 that lower-cases B2 ("Quuk" -> "quuk"),
# and leave the selection on B3 as in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C2").Formula = "=LOWER(B2)"

$ws.Range("B3").Select() | Out-Null
